$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.423852333333334
$ws.Range("H2").Value = 28.271557
$ws.Range("I2").Value = 0.06654336290212845
$ws.Range("J2").Value = 0.06654336290212845
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.083188
$ws.Range("N2").Value = 6.249564
$ws.Range("O2").Value = 0.01853451022101116
$ws.Range("P2").Value = 0.01853451022101116
$ws.Range("Q2").Value = 19.631656094572
$ws.Range("R2").Value = 176.684904851148
$ws.Range("S2").Value = 0.001233348639849955
$ws.Range("T2").Value = 0.001233348639849955

$ws.Range("G3").Value = 9.423852333333334
$ws.Range("H3").Value = 28.271557
$ws.Range("I3").Value = 0.06654336290212845
$ws.Range("J3").Value = 0.06654336290212845
$ws.Range("O3").Value = 0.7177032719746937
$ws.Range("P3").Value = 0.717703271974694
$ws.Range("Q3").Value = 760.1875444965274
$ws.Range("R3").Value = 6841.687900468746
$ws.Range("S3").Value = 0.04775838928305704
$ws.Range("T3").Value = 0.04775838928305706

$ws.Range("G4").Value = 9.423852333333334
$ws.Range("H4").Value = 28.271557
$ws.Range("I4").Value = 0.06654336290212845
$ws.Range("J4").Value = 0.06654336290212845
$ws.Range("M4").Value = 29.09185666666666
$ws.Range("N4").Value = 87.27556999999999
$ws.Range("O4").Value = 0.258835647448298
$ws.Range("P4").Value = 0.258835647448298
$ws.Range("Q4").Value = 274.1573613291656
$ws.Range("R4").Value = 2467.41625196249
$ws.Range("S4").Value = 0.01722379442015947
$ws.Range("T4").Value = 0.01722379442015947

$ws.Range("G5").Value = 9.423852333333334
$ws.Range("H5").Value = 28.271557
$ws.Range("I5").Value = 0.06654336290212845
$ws.Range("J5").Value = 0.06654336290212845
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5537223333333333
$ws.Range("N5").Value = 1.661167
$ws.Range("O5").Value = 0.004926570355997066
$ws.Range("P5").Value = 0.004926570355997067
$ws.Range("Q5").Value = 5.218197503002112
$ws.Range("R5").Value = 46.963777527019
$ws.Range("S5").Value = 0.0003278305590619809
$ws.Range("T5").Value = 0.000327830559061981

$ws.Range("I6").Value = 0.3572423751649123
$ws.Range("J6").Value = 0.3572423751649123
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.083188
$ws.Range("N6").Value = 6.249564
$ws.Range("O6").Value = 0.01853451022101116
$ws.Range("P6").Value = 0.01853451022101116
$ws.Range("Q6").Value = 105.393823602824
$ws.Range("R6").Value = 948.544412425416
$ws.Range("S6").Value = 0.006621312453872371
$ws.Range("T6").Value = 0.006621312453872372

$ws.Range("I7").Value = 0.3572423751649123
$ws.Range("J7").Value = 0.3572423751649123
$ws.Range("O7").Value = 0.7177032719746937
$ws.Range("P7").Value = 0.717703271974694
$ws.Range("S7").Value = 0.2563940215438686
$ws.Range("T7").Value = 0.2563940215438687

$ws.Range("I8").Value = 0.3572423751649123
$ws.Range("J8").Value = 0.3572423751649123
$ws.Range("M8").Value = 29.09185666666666
$ws.Range("N8").Value = 87.27556999999999
$ws.Range("O8").Value = 0.258835647448298
$ws.Range("P8").Value = 0.258835647448298
$ws.Range("Q8").Value = 1471.831639681731
$ws.Range("R8").Value = 13246.48475713558
$ws.Range("S8").Value = 0.09246706147177784
$ws.Range("T8").Value = 0.09246706147177786

$ws.Range("I9").Value = 0.3572423751649123
$ws.Range("J9").Value = 0.3572423751649123
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5537223333333333
$ws.Range("N9").Value = 1.661167
$ws.Range("O9").Value = 0.004926570355997066
$ws.Range("P9").Value = 0.004926570355997067
$ws.Range("Q9").Value = 28.01423295654422
$ws.Range("R9").Value = 252.128096608898
$ws.Range("S9").Value = 0.00175997969539344
$ws.Range("T9").Value = 0.00175997969539344

$ws.Range("G10").Value = 26.84076266666667
$ws.Range("H10").Value = 80.522288
$ws.Range("I10").Value = 0.1895270158659356
$ws.Range("J10").Value = 0.1895270158659356
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.083188
$ws.Range("N10").Value = 6.249564
$ws.Range("O10").Value = 0.01853451022101116
$ws.Range("P10").Value = 0.01853451022101116
$ws.Range("Q10").Value = 55.91435469804801
$ws.Range("R10").Value = 503.2291922824321
$ws.Range("S10").Value = 0.003512790412724927
$ws.Range("T10").Value = 0.003512790412724928

$ws.Range("G11").Value = 26.84076266666667
$ws.Range("H11").Value = 80.522288
$ws.Range("I11").Value = 0.1895270158659356
$ws.Range("J11").Value = 0.1895270158659356
$ws.Range("O11").Value = 0.7177032719746937
$ws.Range("P11").Value = 0.717703271974694
$ws.Range("Q11").Value = 2165.145711357963
$ws.Range("R11").Value = 19486.31140222167
$ws.Range("S11").Value = 0.1360241594145817
$ws.Range("T11").Value = 0.1360241594145817

$ws.Range("G12").Value = 26.84076266666667
$ws.Range("H12").Value = 80.522288
$ws.Range("I12").Value = 0.1895270158659356
$ws.Range("J12").Value = 0.1895270158659356
$ws.Range("M12").Value = 29.09185666666666
$ws.Range("N12").Value = 87.27556999999999
$ws.Range("O12").Value = 0.258835647448298
$ws.Range("P12").Value = 0.258835647448298
$ws.Range("Q12").Value = 780.8476203226843
$ws.Range("R12").Value = 7027.628582904159
$ws.Range("S12").Value = 0.04905634786060327
$ws.Range("T12").Value = 0.04905634786060329

$ws.Range("G13").Value = 26.84076266666667
$ws.Range("H13").Value = 80.522288
$ws.Range("I13").Value = 0.1895270158659356
$ws.Range("J13").Value = 0.1895270158659356
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.5537223333333333
$ws.Range("N13").Value = 1.661167
$ws.Range("O13").Value = 0.004926570355997066
$ws.Range("P13").Value = 0.004926570355997067
$ws.Range("Q13").Value = 14.86232973223289
$ws.Range("R13").Value = 133.760967590096
$ws.Range("S13").Value = 0.0009337181780257037
$ws.Range("T13").Value = 0.000933718178025704

$ws.Range("G14").Value = 54.762539
$ws.Range("H14").Value = 164.287617
$ws.Range("I14").Value = 0.3866872460670236
$ws.Range("J14").Value = 0.3866872460670236
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.083188
$ws.Range("N14").Value = 6.249564
$ws.Range("O14").Value = 0.01853451022101116
$ws.Range("P14").Value = 0.01853451022101116
$ws.Range("Q14").Value = 114.080664094332
$ws.Range("R14").Value = 1026.725976848988
$ws.Range("S14").Value = 0.007167058714563908
$ws.Range("T14").Value = 0.007167058714563909

$ws.Range("G15").Value = 54.762539
$ws.Range("H15").Value = 164.287617
$ws.Range("I15").Value = 0.3866872460670236
$ws.Range("J15").Value = 0.3866872460670236
$ws.Range("O15").Value = 0.7177032719746937
$ws.Range("P15").Value = 0.717703271974694
$ws.Range("Q15").Value = 4417.492823561714
$ws.Range("R15").Value = 39757.43541205542
$ws.Range("S15").Value = 0.2775267017331864
$ws.Range("T15").Value = 0.2775267017331865

$ws.Range("G16").Value = 54.762539
$ws.Range("H16").Value = 164.287617
$ws.Range("I16").Value = 0.3866872460670236
$ws.Range("J16").Value = 0.3866872460670236
$ws.Range("M16").Value = 29.09185666666666
$ws.Range("N16").Value = 87.27556999999999
$ws.Range("O16").Value = 0.258835647448298
$ws.Range("P16").Value = 0.258835647448298
$ws.Range("Q16").Value = 1593.143935290743
$ws.Range("R16").Value = 14338.29541761669
$ws.Range("S16").Value = 0.1000884436957574
$ws.Range("T16").Value = 0.1000884436957574

$ws.Range("G17").Value = 54.762539
$ws.Range("H17").Value = 164.287617
$ws.Range("I17").Value = 0.3866872460670236
$ws.Range("J17").Value = 0.3866872460670236
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.5537223333333333
$ws.Range("N17").Value = 1.661167
$ws.Range("O17").Value = 0.004926570355997066
$ws.Range("P17").Value = 0.004926570355997067
$ws.Range("Q17").Value = 30.32324087433766
$ws.Range("R17").Value = 272.9091678690389
$ws.Range("S17").Value = 0.001905041923515942
$ws.Range("T17").Value = 0.001905041923515942
